$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the added columns
$ws.Range("C1").Value = "fecha_inicio"
$ws.Range("D1").Value = "fecha_final"

# fecha_inicio / fecha_final per periodo row (semester date ranges)
$data = @(
    @{ Row = 2;  Start = "2018-01-01"; End = "2018-06-30" },
    @{ Row = 3;  Start = "2018-07-01"; End = "2018-12-31" },
    @{ Row = 4;  Start = "2019-01-01"; End = "2019-06-30" },
    @{ Row = 5;  Start = "2019-07-01"; End = "2019-12-31" },
    @{ Row = 6;  Start = "2020-01-01"; End = "2020-06-30" },
    @{ Row = 7;  Start = "2020-07-01"; End = "2020-12-31" },
    @{ Row = 8;  Start = "2021-01-01"; End = "2021-06-30" },
    @{ Row = 9;  Start = "2021-07-01"; End = "2021-12-31" },
    @{ Row = 10; Start = "2022-01-01"; End = "2022-06-30" },
    @{ Row = 11; Start = "2022-07-01"; End = "2022-12-31" },
    @{ Row = 12; Start = "2023-01-01"; End = "2023-06-30" },
    @{ Row = 13; Start = "2023-07-01"; End = "2023-12-31" }
)

# Seed the first data cell with a date value, give it the short-date
# number format, then fan that single style out to the whole C:D block
# via copy/paste-format so every cell shares one cellXfs entry (instead
# of minting a fresh style record per cell).
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Value = $data[0].Start
$ws.Range("C2").Copy()
$ws.Range("C2:D13").PasteSpecial(-4122)  # xlPasteFormats

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.Start
    $ws.Cells.Item($r, 4).Value = $item.End
}

$ws.Application.CutCopyMode = $false
[void]$ws.Range("I8").Select()
